# Add the new "Sendikalar" worksheet at the end of the workbook, populate its
# header row, set column width, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add()
$newSheet.Move($wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Sendikalar"

# Header values (reuse existing shared strings: "MusteriId" and "Ad")
$newSheet.Range("A1").Value = "MusteriId"
$newSheet.Range("B1").Value = "Ad"

# Column A width
$newSheet.Columns.Item(1).ColumnWidth = 9.81640625

# Select the full first column, as in the authored file
$newSheet.Range("A1:A1048576").Select()

# Make this new sheet the active tab
$newSheet.Activate()

# Scroll so that tab 2 ("GercekMusteriler", index 1) is the first visible tab
$wb.Windows.Item(1).ScrollWorkbookTabs(1)
